$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string]) {
            if ($val -match '^\[.*\]\(.*\)$' -and $val -notmatch '\{target="_blank"\}$') {
                $cell.Value2 = $val + '{target="_blank"}'
            }
        }
    }
}

$ws.Range("F10").Select()
